# Updates cryptos list values (prices and 1h volume %) per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.037.00"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.518.17"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'608.40"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'148.00"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("D7").Value = "3.517.60"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "'7.93"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "4.112.55"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "3.517.62"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "67.047.07"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'10.70"
$ws.Range("E19").Value = "  +9.01%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "'15.31"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "'437.94"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'0.610"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "3.653.86"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "  -4.09%  "
$ws.Range("D28").Value = "'9.79"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "'8.29"
$ws.Range("E29").Value = "  -4.70%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "'25.45"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "3.510.52"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").Value = "'5.93"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'8.04"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'173.31"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0893"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'5.45"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  -9.96%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'46.02"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'28.17"
$ws.Range("E47").Value = "  -6.67%  "
$ws.Range("D48").Value = "'1.30"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'7.48"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("D51").Value = "'0.990"
$ws.Range("E51").Value = "  +0.31%  "